# Commit: "Fruta / hortaliza, semanal"
# Insert a new weekly record row at row 390 (Vega Monumental Concepción -
# Zapallo / Camote-Paine), pushing the existing rows 390-456 down to
# 391-457, and fill in the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 390..456 down by inserting a fresh row at 390.
$ws.Rows.Item(390).Insert()

# Populate the newly inserted row 390 with the new observation.
$ws.Cells.Item(390, 1).Value2  = 11
$ws.Cells.Item(390, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(390, 3).Value2  = "Bíobío"
$ws.Cells.Item(390, 4).Value2  = 45209
$ws.Cells.Item(390, 5).Value2  = 8
$ws.Cells.Item(390, 6).Value2  = 100112045
$ws.Cells.Item(390, 7).Value2  = "Zapallo"
$ws.Cells.Item(390, 8).Value2  = "Paine"
$ws.Cells.Item(390, 9).Value2  = "1a (guarda)"
$ws.Cells.Item(390, 10).Value2 = 600
$ws.Cells.Item(390, 11).Value2 = 500
$ws.Cells.Item(390, 12).Value2 = 500
$ws.Cells.Item(390, 13).Value2 = 500
$ws.Cells.Item(390, 14).Value2 = "$/kilo (volumen en unidades)"
$ws.Cells.Item(390, 15).Value2 = "Región del Maule"
$ws.Cells.Item(390, 16).Value2 = 500
$ws.Cells.Item(390, 17).Value2 = 1
$ws.Cells.Item(390, 18).Value2 = "Hortaliza"
